$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 32143.777
$ws.Range("I20").Value = 9839.6
$ws.Range("J20").Value = 60024
$ws.Range("K20").Value = 9839.6
$ws.Range("L20").Value = 60024
$ws.Range("M20").Value = -9609.6
$ws.Range("N20").Value = -60484

$ws.Range("H35").Value = 32143.777
$ws.Range("I35").Value = 9839.6
$ws.Range("J35").Value = 60024
$ws.Range("K35").Value = 9839.6
$ws.Range("L35").Value = 60024
$ws.Range("M35").Value = -9460.6
$ws.Range("N35").Value = -60782

$ws.Range("H43").Value = 2343.5557
$ws.Range("J43").Value = 3138.4
$ws.Range("L43").Value = 3138.4
$ws.Range("N43").Value = -3276.4

$ws.Range("H87").Value = 29366.7
$ws.Range("J87").Value = 29366.7
$ws.Range("L87").Value = 29366.7
$ws.Range("N87").Value = -31862.7

$ws.Range("H90").Value = 29366.7
$ws.Range("J90").Value = 29366.7
$ws.Range("L90").Value = 88100.10000000001
$ws.Range("N90").Value = -100580.1

$ws.Range("H111").Value = 991.9091
$ws.Range("I111").Value = 981.1
$ws.Range("K111").Value = 2943.3
$ws.Range("M111").Value = 123.6999999999998

$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

$ws.Range("H132").Value = 3848561
$ws.Range("I132").Value = 4168830.8
$ws.Range("J132").Value = 5325
$ws.Range("K132").Value = 12506492.4
$ws.Range("L132").Value = 15975
$ws.Range("M132").Value = -12503962.4
$ws.Range("N132").Value = -21035

$ws.Range("H141").Value = 405850.34
$ws.Range("I141").Value = 1115.6364
$ws.Range("J141").Value = 1215319.8
$ws.Range("K141").Value = 3346.9092
$ws.Range("L141").Value = 3645959.4
$ws.Range("M141").Value = 1833.0908
$ws.Range("N141").Value = -3656319.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16854
$ws.Range("I32").Value = 14680.86
$ws.Range("J32").Value = 21772.158
$ws.Range("K32").Value = 14680.86
$ws.Range("L32").Value = 21772.158
$ws.Range("M32").Value = -14393.86
$ws.Range("N32").Value = -22346.158

$ws.Range("H61").Value = 2134.6086
$ws.Range("I61").Value = 1247.9459
$ws.Range("J61").Value = 5779.778
$ws.Range("K61").Value = 1247.9459
$ws.Range("L61").Value = 5779.778
$ws.Range("M61").Value = -1035.9459
$ws.Range("N61").Value = -6203.778

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

$ws.Range("H136").Value = 2134.6086
$ws.Range("I136").Value = 1247.9459
$ws.Range("J136").Value = 5779.778
$ws.Range("K136").Value = 3743.8377
$ws.Range("L136").Value = 17339.334
$ws.Range("M136").Value = -1193.8377
$ws.Range("N136").Value = -22439.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 2007.2
$ws.Range("J17").Value = 2007.2
$ws.Range("L17").Value = 2007.2
$ws.Range("N17").Value = -2351.2

$ws.Range("H134").Value = 2314.4736
$ws.Range("I134").Value = 2212.4695
$ws.Range("J134").Value = 2939.25
$ws.Range("K134").Value = 6637.4085
$ws.Range("L134").Value = 8817.75
$ws.Range("M134").Value = -4102.4085
$ws.Range("N134").Value = -13887.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3245.8386
$ws.Range("I31").Value = 1847
$ws.Range("K31").Value = 1847
$ws.Range("M31").Value = -1552

$ws.Range("H34").Value = 3245.8386
$ws.Range("I34").Value = 1847
$ws.Range("K34").Value = 1847
$ws.Range("M34").Value = -1645

$ws.Range("H58").Value = 8199666.5
$ws.Range("I58").Value = 1436.8108
$ws.Range("J58").Value = 20838604
$ws.Range("K58").Value = 1436.8108
$ws.Range("L58").Value = 20838604
$ws.Range("M58").Value = -1233.8108
$ws.Range("N58").Value = -20839010

$ws.Range("H86").Value = 3333
$ws.Range("I86").Value = 1959.8
$ws.Range("K86").Value = 1959.8
$ws.Range("M86").Value = -836.8

$ws.Range("H89").Value = 3333
$ws.Range("I89").Value = 1959.8
$ws.Range("K89").Value = 9799
$ws.Range("M89").Value = -4183

$ws.Range("H95").Value = 23500
$ws.Range("J95").Value = 23500
$ws.Range("L95").Value = 23500
$ws.Range("N95").Value = -28992

$ws.Range("H107").Value = 1160.4474
$ws.Range("I107").Value = 989.5517
$ws.Range("K107").Value = 989.5517
$ws.Range("M107").Value = 930.4483

$ws.Range("H136").Value = 8199666.5
$ws.Range("I136").Value = 1436.8108
$ws.Range("J136").Value = 20838604
$ws.Range("K136").Value = 4310.4324
$ws.Range("L136").Value = 62515812
$ws.Range("M136").Value = -1760.4324
$ws.Range("N136").Value = -62520912

$ws.Range("H141").Value = 31051.725
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 31051.725
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value = 31051.725
$ws.Range("N141").Value = -41411.725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10332.2
$ws.Range("I11").Value = 185
$ws.Range("J11").Value = 17097
$ws.Range("K11").Value = 555
$ws.Range("L11").Value = 51291
$ws.Range("M11").Value = -415
$ws.Range("N11").Value = -51571

$ws.Range("H18").Value = 854.0833
$ws.Range("I18").Value = 425
$ws.Range("J18").Value = 1283.1666
$ws.Range("K18").Value = 1275
$ws.Range("L18").Value = 3849.4998
$ws.Range("M18").Value = -1106
$ws.Range("N18").Value = -4187.4998

$ws.Range("H86").Value = 202
$ws.Range("I86").Value = 202
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 606
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("M86").Value = 580

$ws.Range("H89").Value = 202
$ws.Range("I89").Value = 202
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1818
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("M89").Value = 4110

$ws.Range("H131").Value = 1368.5186
$ws.Range("J131").Value = 1051.5625
$ws.Range("L131").Value = 3154.6875
$ws.Range("N131").Value = -13234.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5950
$ws.Range("I113").Value = 4600
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 4600
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -2430
$ws.Range("N113").Value = -14340

$ws.Range("H126").Value = 559263.5600000001
$ws.Range("I126").Value = 1955
$ws.Range("J126").Value = 837917.8
$ws.Range("K126").Value = 5865
$ws.Range("L126").Value = 2513753.4
$ws.Range("M126").Value = -3395
$ws.Range("N126").Value = -2518693.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 111115400
$ws.Range("I61").Value = 200000750
$ws.Range("J61").Value = 8700
$ws.Range("K61").Value = 200000750
$ws.Range("L61").Value = 8700
$ws.Range("M61").Value = -200000548
$ws.Range("N61").Value = -9104

$ws.Range("H106").Value = 25081.428
$ws.Range("J106").Value = 25081.428
$ws.Range("L106").Value = 25081.428
$ws.Range("N106").Value = -27605.428

$ws.Range("H113").Value = 111115400
$ws.Range("I113").Value = 200000750
$ws.Range("J113").Value = 8700
$ws.Range("K113").Value = 200000750
$ws.Range("L113").Value = 8700
$ws.Range("M113").Value = -199998580
$ws.Range("N113").Value = -13040

$ws.Range("H132").Value = 2644.8667
$ws.Range("I132").Value = 2148.9312
$ws.Range("J132").Value = 3543.75
$ws.Range("K132").Value = 6446.7936
$ws.Range("L132").Value = 10631.25
$ws.Range("M132").Value = -3916.7936
$ws.Range("N132").Value = -15691.25

$ws.Range("H136").Value = 3693
$ws.Range("I136").Value = 2864.3225
$ws.Range("K136").Value = 8592.967500000001
$ws.Range("M136").Value = -6042.967500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1787.5
$ws.Range("I107").Value = 383.33334
$ws.Range("J107").Value = 2630
$ws.Range("K107").Value = 1150.00002
$ws.Range("L107").Value = 7890
$ws.Range("M107").Value = 769.9999800000001
$ws.Range("N107").Value = -11730

$ws.Range("H113").Value = 2670
$ws.Range("I113").Value = 190
$ws.Range("K113").Value = 570
$ws.Range("M113").Value = 1600

$ws.Range("H132").Value = 16872.5
$ws.Range("I132").Value = 4508.48
$ws.Range("J132").Value = 61029.715
$ws.Range("K132").Value = 13525.44
$ws.Range("L132").Value = 183089.145
$ws.Range("M132").Value = -10995.44
$ws.Range("N132").Value = -188149.145

$ws.Range("H136").Value = 2372.2856
$ws.Range("I136").Value = 1966.1852
$ws.Range("J136").Value = 3103.2666
$ws.Range("K136").Value = 5898.5556
$ws.Range("L136").Value = 9309.799800000001
$ws.Range("M136").Value = -3348.5556
$ws.Range("N136").Value = -14409.7998
